# Auto-generated edit script: weekly rotation of Acelga price data
# Shifts rows 86-235 down by 2 (using pre-edit snapshot), sets new data
# for rows 84-85, and appends 2 new rows (236-237) with the data that
# was displaced from the bottom of the shifted window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  [PSCustomObject]@{ Row=84; D=44477; J=2600; K=$null; L=$null; M=$null; N=$null; P=$null; Q=$null }
  [PSCustomObject]@{ Row=85; D=44477; J=1400; K=$null; L=$null; M=$null; N=$null; P=$null; Q=$null }
  [PSCustomObject]@{ Row=86; D=44424; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=87; D=44424; J=1560; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=88; D=44230; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=89; D=44230; J=1700; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=90; D=44382; J=3000; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=91; D=44382; J=1520; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=92; D=44216; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=93; D=44216; J=1680; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=94; D=44370; J=3400; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=95; D=44370; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=96; D=44284; J=2700; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=97; D=44284; J=1440; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=98; D=44174; J=2400; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=99; D=44174; J=1500; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=100; D=44244; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=101; D=44244; J=1700; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=102; D=44449; J=3500; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=103; D=44449; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=104; D=44412; J=3500; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=105; D=44412; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=106; D=44330; J=3460; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=107; D=44330; J=1640; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=108; D=44463; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=109; D=44463; J=1500; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=110; D=44265; J=3600; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=111; D=44265; J=1740; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=112; D=44428; J=3520; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=113; D=44428; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=114; D=44293; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=115; D=44293; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=116; D=44195; J=2700; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=117; D=44195; J=1600; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=118; D=44209; J=3200; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=119; D=44209; J=1700; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=120; D=44447; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=121; D=44447; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=122; D=44421; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=123; D=44421; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=124; D=44232; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=125; D=44232; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=126; D=44316; J=3400; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=127; D=44316; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=128; D=44356; J=3450; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=129; D=44356; J=1660; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=130; D=44169; J=2800; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=131; D=44169; J=1520; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=132; D=44454; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=133; D=44454; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=134; D=44410; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=135; D=44410; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=136; D=44319; J=2800; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=137; D=44319; J=1480; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=138; D=44473; J=2000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=139; D=44473; J=1400; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=140; D=44235; J=3100; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=141; D=44235; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=142; D=44403; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=143; D=44403; J=1560; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=144; D=44291; J=2800; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=145; D=44291; J=1460; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=146; D=44426; J=3500; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=147; D=44426; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=148; D=44351; J=3460; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=149; D=44351; J=1680; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=150; D=44214; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=151; D=44214; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=152; D=44274; J=3300; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=153; D=44274; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=154; D=44445; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=155; D=44445; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=156; D=44344; J=3460; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=157; D=44344; J=1680; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=158; D=44407; J=3560; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=159; D=44407; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=160; D=44358; J=3480; K=500; L=550; M=525; N="`$/atado 1,5 a 2 kilos"; P=262; Q=2 }
  [PSCustomObject]@{ Row=161; D=44358; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=162; D=44391; J=3440; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=163; D=44391; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=164; D=44335; J=3460; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=165; D=44335; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=166; D=44165; J=2600; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=167; D=44165; J=1600; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=168; D=44258; J=3600; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=169; D=44258; J=1760; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=170; D=44475; J=3120; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=171; D=44475; J=1400; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=172; D=44204; J=2800; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=173; D=44204; J=1600; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=174; D=44300; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=175; D=44300; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=176; D=44186; J=2400; K=400; L=500; M=450; N="`$/atado 1 a 1,5 kilos"; P=300; Q=1.5 }
  [PSCustomObject]@{ Row=177; D=44186; J=1600; K=300; L=350; M=325; N="`$/atado 1 a 1,5 kilos"; P=217; Q=1.5 }
  [PSCustomObject]@{ Row=178; D=44349; J=3400; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=179; D=44349; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=180; D=44372; J=3460; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=181; D=44372; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=182; D=44452; J=2900; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=183; D=44452; J=1400; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=184; D=44202; J=2700; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=185; D=44202; J=1600; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=186; D=44435; J=10220; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=187; D=44435; J=4760; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=188; D=44242; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=189; D=44242; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=190; D=44377; J=3400; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=191; D=44377; J=1640; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=192; D=44433; J=3500; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=193; D=44433; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=194; D=44307; J=3400; K=550; L=600; M=575; N="`$/atado 1,5 a 2 kilos"; P=288; Q=2 }
  [PSCustomObject]@{ Row=195; D=44307; J=1600; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=196; D=44468; J=3120; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=197; D=44468; J=1360; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=198; D=44386; J=3460; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=199; D=44386; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=200; D=44384; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=201; D=44384; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=202; D=44263; J=3100; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=203; D=44263; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=204; D=44363; J=3440; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=205; D=44363; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=206; D=44172; J=2700; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=207; D=44172; J=1400; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=208; D=44328; J=3480; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=209; D=44328; J=1640; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=210; D=44321; J=3400; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=211; D=44321; J=1600; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=212; D=44223; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=213; D=44223; J=1680; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=214; D=44298; J=2800; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=215; D=44298; J=1500; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=216; D=44414; J=3520; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=217; D=44414; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=218; D=44333; J=3000; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=219; D=44333; J=1500; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=220; D=44314; J=3600; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=221; D=44314; J=1660; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
  [PSCustomObject]@{ Row=222; D=44466; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=223; D=44466; J=1800; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=224; D=44270; J=2600; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=225; D=44270; J=1400; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=226; D=44438; J=3200; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=227; D=44438; J=1540; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=228; D=44193; J=2800; K=400; L=500; M=450; N="`$/atado 1,5 a 2 kilos"; P=225; Q=2 }
  [PSCustomObject]@{ Row=229; D=44193; J=1560; K=300; L=350; M=325; N="`$/atado 1,5 a 2 kilos"; P=162; Q=2 }
  [PSCustomObject]@{ Row=230; D=44286; J=3400; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=231; D=44286; J=1600; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=232; D=44389; J=3000; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
  [PSCustomObject]@{ Row=233; D=44389; J=1560; K=350; L=400; M=375; N="`$/atado 1,5 a 2 kilos"; P=188; Q=2 }
  [PSCustomObject]@{ Row=234; D=44312; J=3000; K=550; L=600; M=575; N="`$/atado 1,5 a 2 kilos"; P=288; Q=2 }
  [PSCustomObject]@{ Row=235; D=44312; J=1500; K=450; L=500; M=475; N="`$/atado 1,5 a 2 kilos"; P=238; Q=2 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.D -ne $null) { $ws.Cells.Item($r, 4).Value = $u.D }
    if ($u.J -ne $null) { $ws.Cells.Item($r, 10).Value = $u.J }
    if ($u.K -ne $null) { $ws.Cells.Item($r, 11).Value = $u.K }
    if ($u.L -ne $null) { $ws.Cells.Item($r, 12).Value = $u.L }
    if ($u.M -ne $null) { $ws.Cells.Item($r, 13).Value = $u.M }
    if ($u.N -ne $null) { $ws.Cells.Item($r, 14).Value = $u.N }
    if ($u.P -ne $null) { $ws.Cells.Item($r, 16).Value = $u.P }
    if ($u.Q -ne $null) { $ws.Cells.Item($r, 17).Value = $u.Q }
}

# Append the two new rows at the bottom with full column data
$constA = 8
$constB = "Terminal La Palmera de La Serena"
$constC = "Coquimbo"
$constE = 4
$constF = 100112009
$constG = "Acelga"
$constH = "Sin especificar"
$constO = "Provincia del Elquí"
$constR = "Hortaliza"

$newRows = @(
  [PSCustomObject]@{ Row=236; I="Primera"; D=44326; J=3260; K=500; L=600; M=550; N="`$/atado 1,5 a 2 kilos"; P=275; Q=2 }
  [PSCustomObject]@{ Row=237; I="Segunda"; D=44326; J=1500; K=400; L=450; M=425; N="`$/atado 1,5 a 2 kilos"; P=212; Q=2 }
)

foreach ($nr in $newRows) {
    $r = $nr.Row
    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 4).Value = $nr.D
    $ws.Cells.Item($r, 4).NumberFormat = $ws.Range("D84").NumberFormat()
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $nr.I
    $ws.Cells.Item($r, 10).Value = $nr.J
    $ws.Cells.Item($r, 11).Value = $nr.K
    $ws.Cells.Item($r, 12).Value = $nr.L
    $ws.Cells.Item($r, 13).Value = $nr.M
    $ws.Cells.Item($r, 14).Value = $nr.N
    $ws.Cells.Item($r, 15).Value = $constO
    $ws.Cells.Item($r, 16).Value = $nr.P
    $ws.Cells.Item($r, 17).Value = $nr.Q
    $ws.Cells.Item($r, 18).Value = $constR
}
